$d = $word.ActiveDocument
$d.Content.Find.Execute("-Date-02-27-2024-Day-29.docx", $true, $false, $false, $false, $false,
                         $true, 1, $false, "-Date-02-27-2024-Day-28.docx", 2)
